$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two cell texts that actually changed content
$ws.Range("A3").Value = "Omegam_H0 (pour flat LambdaCDM)"
$ws.Range("A4").Value = "Omegam_Omegal (attention indroduit de la courbure !)"

# Update the selected cell in the sheet view
$ws.Range("A3").Select()
